$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("add_interface")
$v = $ws.Range("A1").Value2
Write-Output "A1.Value2=$v"
$t = $ws.Range("A1").Text
Write-Output "A1.Text=$t"
$n = $ws.Name
Write-Output "sheetname=$n"
